$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.833.83"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "3.497.41"
$ws.Range("E3").Value = "  +0.27%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'594.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'172.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.21%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -0.25%  "
$ws.Range("D9").Value = "'0.131"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'7.18"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.97%  "
$ws.Range("D11").Value = "'0.431"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.98%  "
$ws.Range("D12").Value = "4.100.14"
$ws.Range("E12").Value = "  +0.23%  "
$ws.Range("E13").Value = "  -0.23%  "
$ws.Range("D14").Value = "'29.28"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.02%  "
$ws.Range("D15").Value = "66.858.34"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "3.473.33"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("E18").Value = "  -0.09%  "
$ws.Range("D19").Value = "'14.22"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "'394.56"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Value = "'73.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'10.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.78%  "
$ws.Range("E27").Value = "  +0.55%  "
$ws.Range("D28").Value = "'0.997"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("D29").Value = "'6.14"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("D30").Value = "'1.42"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.75%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("D32").Value = "'23.67"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'7.36"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.50%  "
$ws.Range("E34").Value = "  +0.50%  "
$ws.Range("D35").Value = "'162.68"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.53%  "
$ws.Range("D36").Value = "'0.878"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.04%  "
$ws.Range("E37").Value = "  -0.20%  "
$ws.Range("D38").Value = "'6.89"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.62%  "
$ws.Range("D39").Value = "'4.63"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'0.0737"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("B41").Value = "InjectiveProtocol"
$ws.Range("C41").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D41").Value = "'27.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.28%  "
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "2.829.67"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").Value = "'26.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.87%  "
$ws.Range("D44").Value = "'42.68"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.92%  "
$ws.Range("D45").Value = "'2.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.72%  "
$ws.Range("D46").Value = "'0.0302"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.20%  "
$ws.Range("D47").Value = "'337.09"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.21%  "
$ws.Range("D48").Value = "'34.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.08%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "'6.43"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").Value = "'0.840"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.93%  "
